$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Add a new "测试04" section (rows 23-28), built the same way the
# existing "测试02" (rows 9-14) and "测试03" (rows 16-21) sections were:
# a blank separator row, then a header row, then four browser rows.
# ---------------------------------------------------------------------

# Row 23: section header (new shared string "测试04")
$ws.Range("A23").Value = "测试04"

# Row 24: sub-header row, copy formatting from row 10
$ws.Range("B10:C10").Copy()
$ws.Range("B24:C24").PasteSpecial($xlPasteFormats)
$ws.Range("B24:C24").Merge()
$ws.Range("B24").Value = "浏览器正常运行代码"

# Rows 25-28: browser rows, copy formatting from rows 11-14
$srcRows = @(11, 12, 13, 14)
$dstRows = @(25, 26, 27, 28)
$labels = @("火狐", "谷歌chrome", "iPad safari", "微软 Edge")

for ($i = 0; $i -lt 4; $i++) {
    $srcRow = $srcRows[$i]
    $dstRow = $dstRows[$i]

    $ws.Range("A$srcRow").Copy()
    $ws.Range("A$dstRow").PasteSpecial($xlPasteFormats)
    $ws.Range("A$dstRow").Value = $labels[$i]

    $ws.Range("B${srcRow}:C${srcRow}").Copy()
    $ws.Range("B${dstRow}:C${dstRow}").PasteSpecial($xlPasteFormats)
    $ws.Range("B${dstRow}:C${dstRow}").Merge()
    $ws.Range("B$dstRow").Value = "YES"
}

$excel.CutCopyMode = 0

# Update the view to match the new extent (mirrors how the user scrolled
# down to the newly appended section before saving)
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("B28:C28").Select()
